$d = $word.ActiveDocument

# Locate the run of text that needs to be split up (this is the text that
# Word's proofing pass re-segmented into multiple runs separated by
# <w:proofErr/> markers after the font was changed to Calibri).
$rng = $d.Content
$found = $rng.Find.Execute("f(x) = (ax - b)g(x) + R", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text 'f(x) = (ax - b)g(x) + R'"
}

$xml = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'>
<w:body>
<w:p w14:paraId="793D86C8" w14:textId="77777777" w:rsidR="00556D28" w:rsidRPr="00A82893" w:rsidRDefault="00A82893">
<w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Proxima Nova" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00A82893"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Proxima Nova" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>f(x) = (</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Proxima Nova" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>ax</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Proxima Nova" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Proxima Nova" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>b)g</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Proxima Nova" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>(x) + R</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$rng.InsertXML($xml)
